$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48 (shifts existing rows 48..131 down to 49..132)
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new weekly record
$ws.Cells.Item(48, 1).Value = 8
$ws.Cells.Item(48, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(48, 3).Value = "Coquimbo"
$ws.Cells.Item(48, 4).Value = 45070
$ws.Cells.Item(48, 5).Value = 4
$ws.Cells.Item(48, 6).Value = 100114007
$ws.Cells.Item(48, 7).Value = "Jengibre"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 320
$ws.Cells.Item(48, 11).Value = 17000
$ws.Cells.Item(48, 12).Value = 18000
$ws.Cells.Item(48, 13).Value = 17500
$ws.Cells.Item(48, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(48, 15).Value = "Perú"
$ws.Cells.Item(48, 16).Value = 1346
$ws.Cells.Item(48, 17).Value = 13
$ws.Cells.Item(48, 18).Value = "Hortaliza"
